$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column B (rows 2-6)
$ws.Range("B2").Value = 45.05
$ws.Range("B3").Value = 89.68000000000001
$ws.Range("B4").Value = 88.61
$ws.Range("B5").Value = 93.18000000000001
$ws.Range("B6").Value = 94.47

# Add new rows 7-11
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 95

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 95.61

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 95.91

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 96.11

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 96.18000000000001
